# Weekly fruit/vegetable price-list refresh for Nispero, Vega Modelo de Temuco.
# Re-applies the new batch of rows 2-11 (dates, volumes, prices, origins, units)
# over the previous snapshot, cell by cell, matching the committed OOXML diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44511
$ws.Range("M2").Value = 45
$ws.Range("N2").Value = 28000
$ws.Range("O2").Value = 28000
$ws.Range("P2").Value = 28000
$ws.Range("Q2").Value = '$/bandeja 10 kilos'
$ws.Range("R2").Value = 'Provincia de Los Andes'
$ws.Range("S2").Value = 2800
$ws.Range("T2").Value = 10

# Row 3
$ws.Range("D3").Value = 44511
$ws.Range("M3").Value = 45
$ws.Range("N3").Value = 3200
$ws.Range("O3").Value = 3200
$ws.Range("P3").Value = 3200
$ws.Range("Q3").Value = '$/bandeja 10 kilos'
$ws.Range("R3").Value = 'Provincia de Quillota'
$ws.Range("S3").Value = 320
$ws.Range("T3").Value = 10

# Row 4
$ws.Range("D4").Value = 44496
$ws.Range("M4").Value = 55

# Row 5
$ws.Range("D5").Value = 44466
$ws.Range("N5").Value = 11000
$ws.Range("O5").Value = 11000
$ws.Range("P5").Value = 11000
$ws.Range("Q5").Value = '$/bandeja 5 kilos'
$ws.Range("R5").Value = 'La Ligua'
$ws.Range("S5").Value = 2200
$ws.Range("T5").Value = 5

# Row 6
$ws.Range("D6").Value = 44483
$ws.Range("M6").Value = 35
$ws.Range("N6").Value = 10000
$ws.Range("O6").Value = 10000
$ws.Range("P6").Value = 10000
$ws.Range("Q6").Value = '$/bandeja 5 kilos'
$ws.Range("S6").Value = 2000
$ws.Range("T6").Value = 5

# Row 7
$ws.Range("D7").Value = 44503
$ws.Range("M7").Value = 50
$ws.Range("N7").Value = 28000
$ws.Range("O7").Value = 28000
$ws.Range("P7").Value = 28000
$ws.Range("Q7").Value = '$/bandeja 10 kilos'
$ws.Range("S7").Value = 2800
$ws.Range("T7").Value = 10

# Row 8
$ws.Range("D8").Value = 44488
$ws.Range("L8").Value = 'Primera'
$ws.Range("M8").Value = 100
$ws.Range("Q8").Value = '$/bandeja 5 kilos'
$ws.Range("S8").Value = 2400
$ws.Range("T8").Value = 5

# Row 9
$ws.Range("D9").Value = 44166
$ws.Range("L9").Value = 'Segunda'
$ws.Range("M9").Value = 20
$ws.Range("N9").Value = 12000
$ws.Range("O9").Value = 12000
$ws.Range("P9").Value = 12000
$ws.Range("Q9").Value = '$/caja 18 kilos'
$ws.Range("R9").Value = 'La Ligua'
$ws.Range("S9").Value = 667
$ws.Range("T9").Value = 18

# Row 10
$ws.Range("D10").Value = 44519
$ws.Range("M10").Value = 30
$ws.Range("N10").Value = 28000
$ws.Range("O10").Value = 28000
$ws.Range("P10").Value = 28000
$ws.Range("S10").Value = 2800

# Row 11
$ws.Range("D11").Value = 44515
$ws.Range("M11").Value = 80
$ws.Range("R11").Value = 'Provincia de Los Andes'

Write-Host "Applied 69 cell updates."
